# Automatische test-sync: 2025-06-19 10:30:10
# Append the new "Vragen over samenwerking" mail-log row to the Logs sheet,
# extend the conditional-formatting ranges to cover it, and refresh the
# Dashboard's "Overig" category count.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

# New row 6 data (column E / "Antwoord" is intentionally left blank).
$logs.Range("A6").Value = "Vragen over samenwerking"
$logs.Range("B6").Value = "mailmind.test@zohomail.eu"
$logs.Range("C6").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D6").Value = "Overig"
$logs.Range("F6").Value = "2025-06-19 10:28:10"
$logs.Range("G6").Value = "Nee"

# Extend the existing conditional-formatting rules (Categorie / Beantwoord)
# so they keep covering the data range now that it runs through row 6.
$categorieRules = $logs.Range("D2:D5").FormatConditions()
$categorieRules.Item(1).ModifyAppliesToRange($logs.Range("D2:D6"))

$beantwoordRules = $logs.Range("G2:G5").FormatConditions()
$beantwoordRules.Item(1).ModifyAppliesToRange($logs.Range("G2:G6"))

# Dashboard: "Overig" count goes from 2 to 3 with the new row.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 3
